$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column A rows 76-246 to reflect the new Torch module/children enumeration
$ws.Range("A76").Value = 'module 0 VGG('
$ws.Range("A77").Value = '  (features): Sequential('
$ws.Range("A78").Value = '    (0): Conv2d(3, 64, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A79").Value = '    (1): ReLU(inplace=True)'
$ws.Range("A80").Value = '    (2): Conv2d(64, 64, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A81").Value = '    (3): ReLU(inplace=True)'
$ws.Range("A82").Value = '    (4): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A83").Value = '    (5): Conv2d(64, 128, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A84").Value = '    (6): ReLU(inplace=True)'
$ws.Range("A85").Value = '    (7): Conv2d(128, 128, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A86").Value = '    (8): ReLU(inplace=True)'
$ws.Range("A87").Value = '    (9): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A88").Value = '    (10): Conv2d(128, 256, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A89").Value = '    (11): ReLU(inplace=True)'
$ws.Range("A90").Value = '    (12): Conv2d(256, 256, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A91").Value = '    (13): ReLU(inplace=True)'
$ws.Range("A92").Value = '    (14): Conv2d(256, 256, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A93").Value = '    (15): ReLU(inplace=True)'
$ws.Range("A94").Value = '    (16): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A95").Value = '    (17): Conv2d(256, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A96").Value = '    (18): ReLU(inplace=True)'
$ws.Range("A97").Value = '    (19): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A98").Value = '    (20): ReLU(inplace=True)'
$ws.Range("A99").Value = '    (21): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A100").Value = '    (22): ReLU(inplace=True)'
$ws.Range("A101").Value = '    (23): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A102").Value = '    (24): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A103").Value = '    (25): ReLU(inplace=True)'
$ws.Range("A104").Value = '    (26): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A105").Value = '    (27): ReLU(inplace=True)'
$ws.Range("A106").Value = '    (28): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A107").Value = '    (29): ReLU(inplace=True)'
$ws.Range("A108").Value = '    (30): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A109").Value = '  )'
$ws.Range("A110").Value = '  (avgpool): AdaptiveAvgPool2d(output_size=(7, 7))'
$ws.Range("A111").Value = '  (classifier): Sequential('
$ws.Range("A112").Value = '    (0): Linear(in_features=25088, out_features=4096, bias=True)'
$ws.Range("A113").Value = '    (1): ReLU(inplace=True)'
$ws.Range("A114").Value = '    (2): Dropout(p=0.5, inplace=False)'
$ws.Range("A115").Value = '    (3): Linear(in_features=4096, out_features=4096, bias=True)'
$ws.Range("A116").Value = '    (4): ReLU(inplace=True)'
$ws.Range("A117").Value = '    (5): Dropout(p=0.5, inplace=False)'
$ws.Range("A118").Value = '    (6): Linear(in_features=4096, out_features=1000, bias=True)'
$ws.Range("A119").Value = '  )'
$ws.Range("A120").Value = ')'
$ws.Range("A121").Value = 'module 1 Sequential('
$ws.Range("A122").Value = '  (0): Conv2d(3, 64, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A123").Value = '  (1): ReLU(inplace=True)'
$ws.Range("A124").Value = '  (2): Conv2d(64, 64, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A125").Value = '  (3): ReLU(inplace=True)'
$ws.Range("A126").Value = '  (4): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A127").Value = '  (5): Conv2d(64, 128, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A128").Value = '  (6): ReLU(inplace=True)'
$ws.Range("A129").Value = '  (7): Conv2d(128, 128, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A130").Value = '  (8): ReLU(inplace=True)'
$ws.Range("A131").Value = '  (9): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A132").Value = '  (10): Conv2d(128, 256, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A133").Value = '  (11): ReLU(inplace=True)'
$ws.Range("A134").Value = '  (12): Conv2d(256, 256, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A135").Value = '  (13): ReLU(inplace=True)'
$ws.Range("A136").Value = '  (14): Conv2d(256, 256, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A137").Value = '  (15): ReLU(inplace=True)'
$ws.Range("A138").Value = '  (16): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A139").Value = '  (17): Conv2d(256, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A140").Value = '  (18): ReLU(inplace=True)'
$ws.Range("A141").Value = '  (19): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A142").Value = '  (20): ReLU(inplace=True)'
$ws.Range("A143").Value = '  (21): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A144").Value = '  (22): ReLU(inplace=True)'
$ws.Range("A145").Value = '  (23): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A146").Value = '  (24): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A147").Value = '  (25): ReLU(inplace=True)'
$ws.Range("A148").Value = '  (26): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A149").Value = '  (27): ReLU(inplace=True)'
$ws.Range("A150").Value = '  (28): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A151").Value = '  (29): ReLU(inplace=True)'
$ws.Range("A152").Value = '  (30): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A153").Value = ')'
$ws.Range("A154").Value = 'module 2 Conv2d(3, 64, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A155").Value = 'module 3 ReLU(inplace=True)'
$ws.Range("A156").Value = 'module 4 Conv2d(64, 64, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A157").Value = 'module 5 ReLU(inplace=True)'
$ws.Range("A158").Value = 'module 6 MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A159").Value = 'module 7 Conv2d(64, 128, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A160").Value = 'module 8 ReLU(inplace=True)'
$ws.Range("A161").Value = 'module 9 Conv2d(128, 128, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A162").Value = 'module 10 ReLU(inplace=True)'
$ws.Range("A163").Value = 'module 11 MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A164").Value = 'module 12 Conv2d(128, 256, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A165").Value = 'module 13 ReLU(inplace=True)'
$ws.Range("A166").Value = 'module 14 Conv2d(256, 256, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A167").Value = 'module 15 ReLU(inplace=True)'
$ws.Range("A168").Value = 'module 16 Conv2d(256, 256, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A169").Value = 'module 17 ReLU(inplace=True)'
$ws.Range("A170").Value = 'module 18 MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A171").Value = 'module 19 Conv2d(256, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A172").Value = 'module 20 ReLU(inplace=True)'
$ws.Range("A173").Value = 'module 21 Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A174").Value = 'module 22 ReLU(inplace=True)'
$ws.Range("A175").Value = 'module 23 Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A176").Value = 'module 24 ReLU(inplace=True)'
$ws.Range("A177").Value = 'module 25 MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A178").Value = 'module 26 Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A179").Value = 'module 27 ReLU(inplace=True)'
$ws.Range("A180").Value = 'module 28 Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A181").Value = 'module 29 ReLU(inplace=True)'
$ws.Range("A182").Value = 'module 30 Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A183").Value = 'module 31 ReLU(inplace=True)'
$ws.Range("A184").Value = 'module 32 MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A185").Value = 'module 33 AdaptiveAvgPool2d(output_size=(7, 7))'
$ws.Range("A186").Value = 'module 34 Sequential('
$ws.Range("A187").Value = '  (0): Linear(in_features=25088, out_features=4096, bias=True)'
$ws.Range("A188").Value = '  (1): ReLU(inplace=True)'
$ws.Range("A189").Value = '  (2): Dropout(p=0.5, inplace=False)'
$ws.Range("A190").Value = '  (3): Linear(in_features=4096, out_features=4096, bias=True)'
$ws.Range("A191").Value = '  (4): ReLU(inplace=True)'
$ws.Range("A192").Value = '  (5): Dropout(p=0.5, inplace=False)'
$ws.Range("A193").Value = '  (6): Linear(in_features=4096, out_features=1000, bias=True)'
$ws.Range("A194").Value = ')'
$ws.Range("A195").Value = 'module 35 Linear(in_features=25088, out_features=4096, bias=True)'
$ws.Range("A196").Value = 'module 36 ReLU(inplace=True)'
$ws.Range("A197").Value = 'module 37 Dropout(p=0.5, inplace=False)'
$ws.Range("A198").Value = 'module 38 Linear(in_features=4096, out_features=4096, bias=True)'
$ws.Range("A199").Value = 'module 39 ReLU(inplace=True)'
$ws.Range("A200").Value = 'module 40 Dropout(p=0.5, inplace=False)'
$ws.Range("A201").Value = 'module 41 Linear(in_features=4096, out_features=1000, bias=True)'
$ws.Range("A203").Value = 'Children:'
$ws.Range("A204").Value = 'child 0 Sequential('
$ws.Range("A205").Value = '  (0): Conv2d(3, 64, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A206").Value = '  (1): ReLU(inplace=True)'
$ws.Range("A207").Value = '  (2): Conv2d(64, 64, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A208").Value = '  (3): ReLU(inplace=True)'
$ws.Range("A209").Value = '  (4): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A210").Value = '  (5): Conv2d(64, 128, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A211").Value = '  (6): ReLU(inplace=True)'
$ws.Range("A212").Value = '  (7): Conv2d(128, 128, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A213").Value = '  (8): ReLU(inplace=True)'
$ws.Range("A214").Value = '  (9): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A215").Value = '  (10): Conv2d(128, 256, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A216").Value = '  (11): ReLU(inplace=True)'
$ws.Range("A217").Value = '  (12): Conv2d(256, 256, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A218").Value = '  (13): ReLU(inplace=True)'
$ws.Range("A219").Value = '  (14): Conv2d(256, 256, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A220").Value = '  (15): ReLU(inplace=True)'
$ws.Range("A221").Value = '  (16): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A222").Value = '  (17): Conv2d(256, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A223").Value = '  (18): ReLU(inplace=True)'
$ws.Range("A224").Value = '  (19): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A225").Value = '  (20): ReLU(inplace=True)'
$ws.Range("A226").Value = '  (21): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A227").Value = '  (22): ReLU(inplace=True)'
$ws.Range("A228").Value = '  (23): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A229").Value = '  (24): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A230").Value = '  (25): ReLU(inplace=True)'
$ws.Range("A231").Value = '  (26): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A232").Value = '  (27): ReLU(inplace=True)'
$ws.Range("A233").Value = '  (28): Conv2d(512, 512, kernel_size=(3, 3), stride=(1, 1), padding=(1, 1))'
$ws.Range("A234").Value = '  (29): ReLU(inplace=True)'
$ws.Range("A235").Value = '  (30): MaxPool2d(kernel_size=2, stride=2, padding=0, dilation=1, ceil_mode=False)'
$ws.Range("A236").Value = ')'
$ws.Range("A237").Value = 'child 1 AdaptiveAvgPool2d(output_size=(7, 7))'
$ws.Range("A238").Value = 'child 2 Sequential('
$ws.Range("A239").Value = '  (0): Linear(in_features=25088, out_features=4096, bias=True)'
$ws.Range("A240").Value = '  (1): ReLU(inplace=True)'
$ws.Range("A241").Value = '  (2): Dropout(p=0.5, inplace=False)'
$ws.Range("A242").Value = '  (3): Linear(in_features=4096, out_features=4096, bias=True)'
$ws.Range("A243").Value = '  (4): ReLU(inplace=True)'
$ws.Range("A244").Value = '  (5): Dropout(p=0.5, inplace=False)'
$ws.Range("A245").Value = '  (6): Linear(in_features=4096, out_features=1000, bias=True)'
$ws.Range("A246").Value = ')'

# Update the view selection to match the saved workbook state
$ws.Activate()
$ws.Range("A17").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("P49").Select()
